$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Analysis" column entries to lowercase wording
$ws.Range("E2").Value = "metabarcoding"
$ws.Range("E3").Value = "metabarcoding, metagenomics"

# Update selection to a single active cell, matching the saved view state
$ws.Range("E4").Select()
